# Insert four new "PSFBeads*" sheets right before the existing
# "FieldIlluminationDataset" sheet, mirroring the layout used by the other
# dataset/input/output/key-values sheet groups already in the workbook.

$wb = $excel.ActiveWorkbook

$anchor = $wb.Worksheets.Item("FieldIlluminationDataset")

# Each Add($anchor) inserts the new sheet immediately before $anchor, i.e.
# immediately after whatever was most recently inserted there. To end up
# with the sheets in Dataset, Input, Output, KeyMeasurements order (reading
# left-to-right, right before FieldIlluminationDataset) they must be
# created in the reverse order.

# Match the outline/page-margin defaults used by the sibling sheets
# (sheetPr/outlinePr + 0.75"/0.75"/1"/1"/0.5"/0.5" page margins) on every
# new sheet.
function Set-DefaultSheetLayout($sheet) {
    $sheet.Outline.SummaryRow = 1
    $sheet.Outline.SummaryColumn = 1
    $sheet.PageSetup.LeftMargin = 54
    $sheet.PageSetup.RightMargin = 54
    $sheet.PageSetup.TopMargin = 72
    $sheet.PageSetup.BottomMargin = 72
    $sheet.PageSetup.HeaderMargin = 36
    $sheet.PageSetup.FooterMargin = 36
}

# --- PSFBeadsKeyMeasurements -----------------------------------------------
$wsKeyMeasurements = $wb.Worksheets.Add($anchor)
$wsKeyMeasurements.Name = "PSFBeadsKeyMeasurements"
Set-DefaultSheetLayout($wsKeyMeasurements)
$keyMeasurementHeaders = @(
    "nr_of_beads_analyzed",
    "nr_of_beads_discarded_lateral_edge",
    "nr_of_beads_discarded_axial_edge",
    "nr_of_beads_discarded_self_proximity",
    "nr_of_beads_discarded_cluster",
    "nr_of_beads_discarded_fit_quality",
    "fit_quality_z_mean",
    "fit_quality_z_median",
    "fit_quality_z_stdev",
    "fit_quality_y_mean",
    "fit_quality_y_median",
    "fit_quality_y_stdev",
    "fit_quality_x_mean",
    "fit_quality_x_median",
    "fit_quality_x_stdev",
    "resolution_mean_fwhm_z_pixels",
    "resolution_median_fwhm_z_pixels",
    "resolution_stdev_fwhm_z_pixels",
    "resolution_mean_fwhm_y_pixels",
    "resolution_median_fwhm_y_pixels",
    "resolution_stdev_fwhm_y_pixels",
    "resolution_mean_fwhm_x_pixels",
    "resolution_median_fwhm_x_pixels",
    "resolution_stdev_fwhm_x_pixels",
    "resolution_mean_fwhm_z_microns",
    "resolution_median_fwhm_z_microns",
    "resolution_stdev_fwhm_z_microns",
    "resolution_mean_fwhm_y_microns",
    "resolution_median_fwhm_y_microns",
    "resolution_stdev_fwhm_y_microns",
    "resolution_mean_fwhm_x_microns",
    "resolution_median_fwhm_x_microns",
    "resolution_stdev_fwhm_x_microns",
    "resolution_mean_fwhm_lateral_asymmetry_ratio",
    "resolution_median_fwhm_lateral_asymmetry_ratio",
    "resolution_stdev_fwhm_lateral_asymmetry_ratio"
)
for ($i = 0; $i -lt $keyMeasurementHeaders.Length; $i++) {
    $wsKeyMeasurements.Cells.Item(1, $i + 1).Value = $keyMeasurementHeaders[$i]
}

# --- PSFBeadsOutput --------------------------------------------------------
$wsOutput = $wb.Worksheets.Add($anchor)
$wsOutput.Name = "PSFBeadsOutput"
Set-DefaultSheetLayout($wsOutput)
$outputHeaders = @("bead_crops", "analyzed_bead_centroids", "discarded_bead_centroids_lateral_edge", "discarded_bead_centroids_axial_edge", "discarded_bead_centroids_self_proximity", "discarded_bead_centroids_cluster", "discarded_bead_centroids_fit_quality", "key_values", "psf_properties", "psf_z_profiles", "psf_y_profiles", "psf_x_profiles")
for ($i = 0; $i -lt $outputHeaders.Length; $i++) {
    $wsOutput.Cells.Item(1, $i + 1).Value = $outputHeaders[$i]
}

# --- PSFBeadsInput -------------------------------------------------------
$wsInput = $wb.Worksheets.Add($anchor)
$wsInput.Name = "PSFBeadsInput"
Set-DefaultSheetLayout($wsInput)
$inputHeaders = @("psf_beads_image", "min_lateral_distance_factor", "sigma_z", "sigma_y", "sigma_x", "snr_threshold")
for ($i = 0; $i -lt $inputHeaders.Length; $i++) {
    $wsInput.Cells.Item(1, $i + 1).Value = $inputHeaders[$i]
}

# --- PSFBeadsDataset ---------------------------------------------------
$wsDataset = $wb.Worksheets.Add($anchor)
$wsDataset.Name = "PSFBeadsDataset"
Set-DefaultSheetLayout($wsDataset)
$datasetHeaders = @("input", "output", "microscope", "sample", "experimenter", "acquisition_datetime", "processed", "processing_datetime", "processing_log", "comment", "name", "description")
for ($i = 0; $i -lt $datasetHeaders.Length; $i++) {
    $wsDataset.Cells.Item(1, $i + 1).Value = $datasetHeaders[$i]
}

# Restore the originally active sheet (inserting sheets shifts focus onto
# the last-created one).
$wb.Worksheets.Item("ArgolightBDataset").Activate()
